$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 gains two new trailing columns (F, G) holding paths to the newly
# added metrics/cluster CSVs, shifting the former contents of B2:E2 one
# column to the right (B2 keeps "neurons", duplicated from A2), and J2
# gains a new "NO" value (rename_features/barcode_suffix column).
#
# Target row 2 layout:
# A2 = neurons
# B2 = neurons   (new - duplicate of A2)
# C2 = 10x
# D2 = RNA
# E2 = datasets/10x_1M_neurons/filtered_feature_bc_matrix.h5
# F2 = datasets/10x_1M_neurons/metrics_summary.csv   (new)
# G2 = datasets/10x_1M_neurons/analysis/kmeans/10_clusters/clusters.csv  (new)
# I2 = ENSEMBL
# J2 = NO  (new)

$ws.Range("B2").Value = "neurons"
$ws.Range("C2").Value = "10x"
$ws.Range("D2").Value = "RNA"
$ws.Range("E2").Value = "datasets/10x_1M_neurons/filtered_feature_bc_matrix.h5"
$ws.Range("F2").Value = "datasets/10x_1M_neurons/metrics_summary.csv"
$ws.Range("G2").Value = "datasets/10x_1M_neurons/analysis/kmeans/10_clusters/clusters.csv"
$ws.Range("I2").Value = "ENSEMBL"
$ws.Range("J2").Value = "NO"

# Move the active-cell selection from I3 to J3, matching the author's
# final cursor position.
$ws.Range("J3").Select()
